$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Wookie Mistakes")

# --- Week-1 start date moved back one day (8/31 -> 8/30) for both tables ---
# Columns C:Q are formulas (=prev+7) and recalc automatically.
$ws.Range("B2").Value = 44803
$ws.Range("B14").Value = 44803

# --- Table 1 (rows 3-10): after-match attendance inputs, replacing "A" placeholders ---
$ws.Range("E3").Value = "DNP"
$ws.Range("I3").Value = "DNP"

$ws.Range("E4").Value = "W"
$ws.Range("I4").Value = "DNP"

$ws.Range("E5").Value = "W"
$ws.Range("I5").Value = "DNP"

$ws.Range("E6").Value = "W"
$ws.Range("I6").Value = "DNP"

$ws.Range("E7").Value = "W"
$ws.Range("I7").Value = "DNP"

$ws.Range("E8").Value = "DNP"
$ws.Range("I8").Value = "DNP"

$ws.Range("E9").Value = "DNP"
$ws.Range("I9").Value = "DNP"

$ws.Range("E10").Value = "L"
$ws.Range("I10").Value = "DNP"

# Games-available count correction for row 3
$ws.Range("S3").Value = 6

# --- Table 2 (rows 15-22): after-match attendance inputs, replacing "A" placeholders ---
$ws.Range("E15").Value = "L"
$ws.Range("E16").Value = "W"
$ws.Range("E17").Value = "W"
$ws.Range("E18").Value = "W"
$ws.Range("E19").Value = "DNP"
$ws.Range("E20").Value = "DNP"
$ws.Range("E21").Value = "DNP"
$ws.Range("E22").Value = "W"

# A22 previously carried a stray "apply fill" style flag; drop it so the
# cell's effective format matches the other rows in its column (border-only).
$ws.Range("A22").Interior.Pattern = -4142

# --- View state: active cell moved to T5 ---
$ws.Range("T5").Select()
